# config.xlsx - "locales" sheet text/typo fixes + selection/row-height update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) shop_item_ship_desc (ru) - C20: capitalize fix "Гибели" -> "гибели"
$ws.Range("C20").Value = "Дополнительный корабль, который`nпозволит продолжить игру после`nгибели"

# --- 2) shop_item_teleport_desc (en) - B22: re-wrap the paragraph
$ws.Range("B22").Value = "Teleporter`nwhen key SPACEBAR pressed`nThen move ship at cursor pos`ninstantly. After jump ship get 3sec`ninvulnerability"

# --- 3) shop_item_teleport_desc (ru) - C22: re-wrap + typo fixes (Точку -> точку, Неузявимость -> неузявимость)
$ws.Range("C22").Value = "Телепорт позволяет по нажатию`nкнопки ПРОБЕЛ переместиться`nв точку положения прицела.`nПосле прижка корабль получает`nнеузявимость на 3 секунды."

# --- 4) shop_item_resonator_desc (en) - B28: re-wrap the paragraph
$ws.Range("B28").Value = "If press SPACEBAR inner 3 seconds`nafter hyper jump appear wave that hit`nasteroids by 1 hp."

# --- 5) Row 20 height grew to fit the rewrapped text (35.05 -> 46.25)
$ws.Rows.Item(20).RowHeight = 46.25

# --- 6) view state: select C21 and scroll the window so row 8 is the top visible row
$ws.Range("C21").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
